# Dia 2 (segundo "Dia"), Hora arreglada:
# "Duración: 1 Hora" -> "Duración: 1 Hora y 30 minutos"
# Only the SECOND occurrence (the "Dia 7 de junio del 2021" section,
# i.e. the paragraph right before the trailing empty paragraph) is
# affected; the first "Dia 6" section's "1 Hora" stays untouched.

$d = $word.ActiveDocument

$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Duración: 1 Hora") {
        $targetParagraph = $p
    }
}

$pEnd = $targetParagraph.Range.End
$r = $d.Range($pEnd - 1, $pEnd - 1)
$r.InsertAfter(" y 30 minutos")

# Match the formatting of the preceding " 1 Hora" run (Times New Roman,
# 12pt / sz 24 half-points, es-ES) on the newly inserted text.
$newRun = $d.Range($pEnd - 1, $pEnd - 1 + 13)
$newRun.Font.Name = "Times New Roman"
$newRun.Font.Size = 12
$newRun.LanguageID = 3082
